# Update column C (Förändrad) dates from 45233 (2023-11-03) to 45243 (2023-11-13)
# for rows 2 through 10, preserving existing cell formatting/style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45233) {
        $cell.Value2 = 45243
    }
}
